$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Make sure a proper "Hyperlink" character style exists (matching the one
#    Word normally mints the first time a hyperlink is inserted) so the new
#    hyperlink run below can reference it.
# ---------------------------------------------------------------------------
$hlStyle = $d.Styles.Add("Hyperlink", 2)      # 2 = wdStyleTypeCharacter
$hlStyle.BaseStyle = "Standardskrifttypeiafsnit"
$hlStyle.Priority = 99
$hlStyle.UnhideWhenUsed = $true
$hlFont = $hlStyle.Font
$hlFont.Color = 12673797                      # RGB 0563C1 (wdColor BGR order)
$hlFont.Underline = 1                          # wdUnderlineSingle

# ---------------------------------------------------------------------------
# 2. Turn the first empty paragraph (paragraph 2) into the new "noted link"
#    paragraph: a hyperlink to the snook.ca article, a trailing space, and
#    the _GoBack bookmark (moved here from the "Tilfoj et sogefelt?" para).
# ---------------------------------------------------------------------------
$linkPara = $d.Paragraphs.Item(2)
$linkRange = $linkPara.Range
$linkRange.Collapse(1)                         # 1 = wdCollapseStart
$url = "http://snook.ca/archives/javascript/simplest-jquery-slideshow"
$hyperlink = $d.Hyperlinks.Add($linkRange, $url)

# Append a literal space right after the hyperlink text, inside the same
# paragraph (use the hyperlink's own range so the insertion point can't
# slip into the following paragraph).
$afterLink = $hyperlink.Range
$afterLink.Collapse(0)                         # 0 = wdCollapseEnd
# Insert a temporary marker character after the space so the insertion
# point used for the bookmark below is never flush against the paragraph
# mark (that edge case makes Bookmarks.Add snap to the whole paragraph).
$afterLink.InsertAfter(" X")

$linkPara = $d.Paragraphs.Item(2)
$markerPos = $linkPara.Range.End - 2

# Relocate the _GoBack bookmark from the "Tilfoj et sogefelt?" paragraph to
# the end of this paragraph (right after the space, before the marker).
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()
$bookmarkRange = $d.Range($markerPos, $markerPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# Remove the temporary marker character now that the bookmark is anchored.
$markerRange = $d.Range($markerPos, $markerPos + 1)
$markerRange.Delete()

# ---------------------------------------------------------------------------
# 3. Mark up the footer paragraph (last paragraph) the way Word's proofer
#    does: spell/grammar-check bracketing around "Footer" / "om" and split
#    "TagXpert" out as its own spell-checked run.
# ---------------------------------------------------------------------------
$footerPara = $d.Paragraphs.Last
$footerRange = $footerPara.Range

$footerXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00A32927" w:rsidRDefault="00A97497">
<w:proofErr w:type="spellStart"/>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>F</w:t></w:r>
<w:r w:rsidR="00A32927"><w:t>ooter</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t>&#8230;</w:t></w:r>
<w:r><w:br/><w:t>om</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t xml:space="preserve"> info om selv, og i kasse ved siden af de logo-links til </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>TagXpert</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> og andet</w:t></w:r>
</w:p>
'@

$footerRange.InsertXML($footerXml)
